# Generate Report for Handback
#
# For each localized-language sheet (zh-cn, de-de) this:
#   - updates the Status column (C) to reflect that the file has been
#     handed back and is in sync with en-US
#   - fills in the "Latest Target File" (F) and "Latest Handback File" (G)
#     columns with hyperlinked file names (mirroring the existing
#     Source File Name / Latest Handoff File hyperlinks)
#   - stamps the "Latest Handback DateTime" (H) column with the handback
#     timestamp

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

$sheets = @(
    @{
        Name = "zh-cn"
        HandbackDateTime = "2016-03-20 22:16:27"
        Row2 = @{
            SourceName = "6db2549d-b57e-451f-b4ef-d2dfe809bde2.md"
            SourceUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/e772fc90a8786ea21d7b5153a0232d85550e5962/e2e/6db2549d-b57e-451f-b4ef-d2dfe809bde2.md"
            XlfName    = "6db2549d-b57e-451f-b4ef-d2dfe809bde2.f68d894b023cc27dc44811613346b78582a35dfc.zh-cn.xlf"
            XlfUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/89749e36bb1725923f0fca4d1002d191ee1e343e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/6db2549d-b57e-451f-b4ef-d2dfe809bde2.f68d894b023cc27dc44811613346b78582a35dfc.zh-cn.xlf"
        }
        Row3 = @{
            SourceName = "b151fd36-158f-4d8b-b214-eebf9c23b2fe.md"
            SourceUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/e772fc90a8786ea21d7b5153a0232d85550e5962/e2e/b151fd36-158f-4d8b-b214-eebf9c23b2fe.md"
            XlfName    = "b151fd36-158f-4d8b-b214-eebf9c23b2fe.4492a079f47be2a90691bd3d03a22eec314412d4.zh-cn.xlf"
            XlfUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/89749e36bb1725923f0fca4d1002d191ee1e343e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/b151fd36-158f-4d8b-b214-eebf9c23b2fe.4492a079f47be2a90691bd3d03a22eec314412d4.zh-cn.xlf"
        }
    },
    @{
        Name = "de-de"
        HandbackDateTime = "2016-03-20 22:16:34"
        Row2 = @{
            SourceName = "6db2549d-b57e-451f-b4ef-d2dfe809bde2.md"
            SourceUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/e772fc90a8786ea21d7b5153a0232d85550e5962/e2e/6db2549d-b57e-451f-b4ef-d2dfe809bde2.md"
            XlfName    = "6db2549d-b57e-451f-b4ef-d2dfe809bde2.f68d894b023cc27dc44811613346b78582a35dfc.de-de.xlf"
            XlfUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/551907b7d00ecbeaf4a835d7aa201e863e120a68/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/6db2549d-b57e-451f-b4ef-d2dfe809bde2.f68d894b023cc27dc44811613346b78582a35dfc.de-de.xlf"
        }
        Row3 = @{
            SourceName = "b151fd36-158f-4d8b-b214-eebf9c23b2fe.md"
            SourceUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/e772fc90a8786ea21d7b5153a0232d85550e5962/e2e/b151fd36-158f-4d8b-b214-eebf9c23b2fe.md"
            XlfName    = "b151fd36-158f-4d8b-b214-eebf9c23b2fe.4492a079f47be2a90691bd3d03a22eec314412d4.de-de.xlf"
            XlfUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/551907b7d00ecbeaf4a835d7aa201e863e120a68/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/b151fd36-158f-4d8b-b214-eebf9c23b2fe.4492a079f47be2a90691bd3d03a22eec314412d4.de-de.xlf"
        }
    }
)

foreach ($sheetInfo in $sheets) {
    $ws = $wb.Worksheets.Item($sheetInfo.Name)

    foreach ($rowNum in @(2, 3)) {
        if ($rowNum -eq 2) { $rowData = $sheetInfo.Row2 } else { $rowData = $sheetInfo.Row3 }

        # Status column
        $ws.Range("C$rowNum").Value = $statusHandedBack

        # Latest Target File column (F) -- hyperlink to the source markdown file
        $fCell = $ws.Range("F$rowNum")
        $fCell.Value = $rowData.SourceName
        $ws.Hyperlinks.Add($fCell, $rowData.SourceUrl, "", "", $rowData.SourceName)
        $fCell.Style = "HyperLink"

        # Latest Handback File column (G) -- hyperlink to the translated xlf file
        $gCell = $ws.Range("G$rowNum")
        $gCell.Value = $rowData.XlfName
        $ws.Hyperlinks.Add($gCell, $rowData.XlfUrl, "", "", $rowData.XlfName)
        $gCell.Style = "HyperLink"

        # Latest Handback DateTime column (H)
        $ws.Range("H$rowNum").Value = $sheetInfo.HandbackDateTime
    }
}
